$d = $word.ActiveDocument
$d.TrackRevisions = $false

# Find the paragraph ending in "... host." (User account name description)
$rng = $d.Content
$found = $rng.Find.Execute(" host.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Move to the very end of this found range (after "host.")
    $insertPoint = $rng.Duplicate
    $insertPoint.Collapse(0)  # wdCollapseEnd

    # Insert the new text after "host."
    $insertPoint.InsertAfter(" A Microsoft account cannot be used with SMBSync2. Please create a local account and use it.")
}
